# "Added min/max on criteria but to manage"
#
# Insert a new row right below the header row that records, for every
# criterion column (Prix / Vitesse / Robustesse), whether it should be
# minimised or maximised. The new row's first cell is the label
# "Min/Max" and the three criterion cells default to "Max". The row
# that used to be row 2 ("Poids" ...) simply shifts down to row 3 with
# its content untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2..8 down to 3..9 and open up a blank row 2.
$ws.Rows.Item(2).Insert()

# Populate the new row.
$ws.Range("A2").Value = "Min/Max"
$ws.Range("B2").Value = "Max"
$ws.Range("C2").Value = "Max"
$ws.Range("D2").Value = "Max"

# Re-apply the sheet's "Normal" style across the whole used range so
# every cell (old and new) shares the same cell format.
$ws.Range("A1:D9").Style = "Normal"

# Restore the cursor to the newly inserted row, mirroring the saved
# selection in the edited workbook.
$null = $ws.Range("A3").Select()
